# P3 Order Files Updated
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated ConditionType (column C) values for the remaining 16 trials
$conditionType = @(1,1,2,4,1,4,3,2,3,4,4,2,2,1,3,3)

# New ITI column (column D) values
$iti = @(9,8,7,7,6,8,7,6,7,6,6,6,6,6,8,9)

# Header for new column D
$ws.Cells.Item(1, 4).Value = "ITI"

for ($i = 0; $i -lt $conditionType.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 3).Value = $conditionType[$i]
    $ws.Cells.Item($row, 4).Value = $iti[$i]
}

# Remove the trailing trials (rows 18-20) that are no longer part of the order
[void]$ws.Rows("18:20").Delete()

# Widen column F to match the updated layout
[void]($ws.Columns("F:F").ColumnWidth = 16)

# Update selection to match the saved view state
[void]$ws.Range("D27").Select()
